# Change the closing date on the "Pagina de Requistos FECHADA" line
# from 13/05/2013 to 16/05/2013.
#
# The committed OOXML shows the edited run split into two runs that share
# identical rPr (<w:sz w:val="40"/><w:szCs w:val="40"/>):
#   " data 16"  +  "/05/2013"
# We reproduce that by editing only the "13" substring and then touching
# the (unchanged) tail run's own formatting, which makes the engine keep
# the tail as its own <w:r> instead of fusing it back into the edited run.

$d = $word.ActiveDocument

# Locate "13/05/2013" in the last paragraph without relying on hard-coded
# absolute character offsets.
$para = $d.Paragraphs($d.Paragraphs.Count).Range
$paraText = $para.Text
$relIdx = $paraText.IndexOf("13/05/2013")
$absStart = $para.Start + $relIdx

# Replace just "13" -> "16"
$oldNum = $d.Range($absStart, $absStart + 2)
$oldNum.Text = "16"

# Touch the (unchanged) tail run's formatting so the engine keeps it as a
# distinct run with its own <w:rPr> instead of coalescing it back into the
# preceding run.
$tail = $d.Range($absStart + 2, $absStart + 2 + 8)
$tailBold = $tail.Bold
$tail.Bold = 1
$tail.Bold = $tailBold
